$d = $word.ActiveDocument

# --- Block 1: rework the "Decision tree" notes section ---
$p1 = $d.Paragraphs.Item(12)
$p2 = $d.Paragraphs.Item(20)
$start = $p1.Range.Start
$end = $p2.Range.End
$r = $d.Range($start, $end)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Decision</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>tre</w:t>
      </w:r>
      <w:r>
        <w:t>e</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> (Variablennamen checken ab </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>supervised</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> learning)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Decision</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>tree</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> verbessern</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:t>Andere Klassifikationsmethoden</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Feature </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>selection</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Maybe irgendwie</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>overfitting</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:t>Stabilität testen</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:t>Visualisierung lernkurven</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:t>Neuronales netz vielleicht</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:ind w:left="1440"/>
      </w:pPr>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Host "block1 done"
